# Updates the cryptos list (prices / 1h volume %) and fixes two pairs of
# rows whose rank order had flipped since the last scrape.
#
# All values in columns B-E are stored as literal text in the workbook
# (e.g. "1.000", "29.263.36" are not real numbers), so every write below
# is done via Range.Formula with a leading apostrophe. That forces Excel
# to keep the text exactly as given (no numeric coercion, no loss of
# trailing zeros, no re-interpretation of the two-dot "thousands" prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    $ws.Range($addr).Formula = "'" + $value
}

$updates = [ordered]@{
    "D2"  = "29.275.25";  "E2"  = "  +0.45%  "
    "D3"  = "1.892.94";   "E3"  = "  -0.40%  "
    "D4"  = "1.001";      "E4"  = "  +0.05%  "
    "D5"  = "322.70";     "E5"  = "  -2.92%  "
    "D6"  = "1.000";      "E6"  = "  +0.03%  "
    "D7"  = "0.4719";     "E7"  = "  +2.42%  "
    "D8"  = "0.4041";     "E8"  = "  -1.89%  "
    "D9"  = "47.35";      "E9"  = "  -1.06%  "
    "D10" = "0.08016";    "E10" = "  -0.01%  "
    "D11" = "0.9962";     "E11" = "  -1.45%  "
    "D12" = "23.07";      "E12" = "  +4.35%  "
    "D13" = "1.973.75";   "E13" = "  +3.79%  "
    "D14" = "5.946";      "E14" = "  +0.05%  "
    "D15" = "7.042";      "E15" = "  -0.97%  "
    "D16" = "89.42";      "E16" = "  +0.25%  "
    "E17" = "  +0.05%  "
    "E18" = "  +1.00%  "
    "D20" = "17.51";      "E20" = "  -0.66%  "
    "D21" = "0.9998";     "E21" = "  -0.03%  "
    "D22" = "29.274.40";  "E22" = "  +0.51%  "
    "E23" = "  -0.04%  "
    "D24" = "11.71";      "E24" = "  +2.75%  "
    "D26" = "2.072.06";   "E26" = "  -2.50%  "
    "D27" = "155.38";     "E27" = "  -1.10%  "
    "D28" = "19.69";      "E28" = "  -0.21%  "
    "D29" = "5.958";      "E29" = "  +5.69%  "
    "D30" = "2.086";      "E30" = "  -1.59%  "
    "D31" = "117.36";     "E31" = "  +0.15%  "
    "D32" = "1.023";      "E32" = "  -1.68%  "
    "D33" = "0.09412";    "E33" = "  -0.01%  "

    # Rows 34/35 swapped places (HuobiToken <-> ARBITRUM) plus new figures.
    "B34" = "ARBITRUM";   "C34" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D34" = "1.384";      "E34" = "  -2.88%  "
    "B35" = "HuobiToken"; "C35" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D35" = "3.528";      "E35" = "  -0.47%  "

    "D36" = "5.360";      "E36" = "  +0.07%  "

    # Rows 37/38 swapped places (Hedera <-> VeChain) plus new figures.
    "B37" = "VeChain";    "C37" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D37" = "0.02245";    "E37" = "  +0.01%  "
    "B38" = "Hedera";     "C38" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D38" = "0.06047";    "E38" = "  -0.85%  "

    "D39" = "1.172";      "E39" = "  -0.54%  "
    "D40" = "7.995";      "E40" = "  -5.16%  "
    "D41" = "0.5828";     "E41" = "  -0.32%  "
    "D42" = "0.1827";     "E42" = "  +0.14%  "
    "E43" = "  -0.89%  "
    "D44" = "1.276";      "E44" = "  +1.54%  "

    # Rows 45/46 swapped places (RenderToken <-> Cronos) plus new figures.
    "B45" = "Cronos";       "C45" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D45" = "0.07703";      "E45" = "  +2.86%  "
    "B46" = "RenderToken";  "C46" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D46" = "2.363";        "E46" = "  +1.14%  "

    "D47" = "12.22";      "E47" = "  +1.05%  "
    "D48" = "0.5480";     "E48" = "  -1.09%  "
    "D49" = "1.909";      "E49" = "  -0.74%  "
    "D50" = "113.24";     "E50" = "  +0.39%  "
    "D51" = "0.2972";     "E51" = "  +1.85%  "
}

foreach ($addr in $updates.Keys) {
    Set-Text $addr $updates[$addr]
}
